$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 727.95
$ws.Range("I33").Value = 838.25
$ws.Range("K33").Value = 838.25
$ws.Range("M33").Value = -609.25
$ws.Range("H88").Value = 18563226
$ws.Range("J88").Value = 53502.2
$ws.Range("L88").Value = 53502.2
$ws.Range("N88").Value = -54314.2
$ws.Range("H91").Value = 18563226
$ws.Range("J91").Value = 53502.2
$ws.Range("L91").Value = 53502.2
$ws.Range("N91").Value = -56310.2
$ws.Range("H137").Value = 3378.3076
$ws.Range("I137").Value = 4131.143
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 12393.429
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -9843.429
$ws.Range("N137").Value = -12600
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2157428.5
$ws.Range("I32").Value = 2234374.8
$ws.Range("K32").Value = 2234374.8
$ws.Range("M32").Value = -2234087.8
$ws.Range("H45").Value = 3987.1765
$ws.Range("I45").Value = 861.4
$ws.Range("J45").Value = 8452.571
$ws.Range("K45").Value = 861.4
$ws.Range("L45").Value = 8452.571
$ws.Range("M45").Value = -484.4
$ws.Range("N45").Value = -9206.571
$ws.Range("H61").Value = 5992.7954
$ws.Range("I61").Value = 2929.0322
$ws.Range("K61").Value = 2929.0322
$ws.Range("M61").Value = -2717.0322
$ws.Range("H74").Value = 37967.535
$ws.Range("I74").Value = 51035.55
$ws.Range("J74").Value = 5297.5
$ws.Range("K74").Value = 51035.55
$ws.Range("L74").Value = 5297.5
$ws.Range("M74").Value = -50161.55
$ws.Range("N74").Value = -7045.5
$ws.Range("H77").Value = 37967.535
$ws.Range("I77").Value = 51035.55
$ws.Range("J77").Value = 5297.5
$ws.Range("K77").Value = 255177.75
$ws.Range("L77").Value = 26487.5
$ws.Range("M77").Value = -250809.75
$ws.Range("N77").Value = -35223.5
$ws.Range("H97").Value = 4909499.5
$ws.Range("I97").Value = 479.9091
$ws.Range("J97").Value = 13909368
$ws.Range("K97").Value = 479.9091
$ws.Range("L97").Value = 13909368
$ws.Range("M97").Value = 16.09089999999998
$ws.Range("N97").Value = -13910360
$ws.Range("H102").Value = 2548.8635
$ws.Range("I102").Value = 2432.1428
$ws.Range("K102").Value = 2432.1428
$ws.Range("M102").Value = -810.1428000000001
$ws.Range("H109").Value = 69321
$ws.Range("J109").Value = 69321
$ws.Range("L109").Value = 69321
$ws.Range("N109").Value = -72095
$ws.Range("H110").Value = 55557576
$ws.Range("I110").Value = 2066.6667
$ws.Range("K110").Value = 2066.6667
$ws.Range("M110").Value = -21.66670000000022
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").ClearContents()
$ws.Range("N125").Value = 0
$ws.Range("H132").Value = 8917.450999999999
$ws.Range("I132").Value = 8463.666999999999
$ws.Range("K132").Value = 25391.001
$ws.Range("M132").Value = -22861.001
$ws.Range("H133").Value = 80000
$ws.Range("J133").Value = 80000
$ws.Range("L133").Value = 80000
$ws.Range("N133").Value = -85060
$ws.Range("H136").Value = 5992.7954
$ws.Range("I136").Value = 2929.0322
$ws.Range("K136").Value = 8787.096600000001
$ws.Range("M136").Value = -6237.096600000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7034.0166
$ws.Range("I31").Value = 3224.0417
$ws.Range("J31").Value = 9574
$ws.Range("K31").Value = 3224.0417
$ws.Range("L31").Value = 9574
$ws.Range("M31").Value = -2929.0417
$ws.Range("N31").Value = -10164
$ws.Range("H34").Value = 7034.0166
$ws.Range("I34").Value = 3224.0417
$ws.Range("J34").Value = 9574
$ws.Range("K34").Value = 3224.0417
$ws.Range("L34").Value = 9574
$ws.Range("M34").Value = -3022.0417
$ws.Range("N34").Value = -9978
$ws.Range("H99").Value = 5871.1665
$ws.Range("I99").Value = 2600
$ws.Range("K99").Value = 2600
$ws.Range("M99").Value = -1102
$ws.Range("H122").Value = 3210.0557
$ws.Range("I122").Value = 3150.2144
$ws.Range("K122").Value = 9450.643199999999
$ws.Range("M122").Value = -7000.643199999999
$ws.Range("H126").Value = 5871.1665
$ws.Range("I126").Value = 2600
$ws.Range("K126").Value = 7800
$ws.Range("M126").Value = -5330
$ws.Range("H132").Value = 5145.837
$ws.Range("I132").Value = 3489.8518
$ws.Range("K132").Value = 10469.5554
$ws.Range("M132").Value = -7939.555399999999
$ws.Range("H134").Value = 6716.645
$ws.Range("I134").Value = 2991.9285
$ws.Range("J134").Value = 9784.058999999999
$ws.Range("K134").Value = 8975.7855
$ws.Range("L134").Value = 29352.177
$ws.Range("M134").Value = -6440.7855
$ws.Range("N134").Value = -34422.177
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 12821434
$ws.Range("J92").Value = 15385641
$ws.Range("L92").Value = 46156923
$ws.Range("N92").Value = -46159419
$ws.Range("H107").Value = 10527075
$ws.Range("I107").Value = 609.8
$ws.Range("J107").Value = 14286527
$ws.Range("K107").Value = 1829.4
$ws.Range("L107").Value = 42859581
$ws.Range("M107").Value = 90.60000000000014
$ws.Range("N107").Value = -42863421
$ws.Range("H140").Value = 161936.69
$ws.Range("I140").Value = 191977
$ws.Range("K140").Value = 575931
$ws.Range("M140").Value = -570751
$ws.Range("H141").Value = 5787.7144
$ws.Range("I141").Value = 2113.818
$ws.Range("J141").Value = 19258.666
$ws.Range("K141").Value = 6341.454000000001
$ws.Range("L141").Value = 57775.99800000001
$ws.Range("M141").Value = -1161.454000000001
$ws.Range("N141").Value = -68135.99800000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8766.200000000001
$ws.Range("I43").Value = 8766.200000000001
$ws.Range("K43").Value = 8766.200000000001
$ws.Range("M43").Value = -8615.200000000001
$ws.Range("H97").Value = 2053.5264
$ws.Range("I97").Value = 1867.3077
$ws.Range("K97").Value = 1867.3077
$ws.Range("M97").Value = -1371.3077
$ws.Range("H113").Value = 7890.1816
$ws.Range("I113").Value = 5031.8335
$ws.Range("J113").Value = 9523.522999999999
$ws.Range("K113").Value = 5031.8335
$ws.Range("L113").Value = 9523.522999999999
$ws.Range("M113").Value = -2861.8335
$ws.Range("N113").Value = -13863.523
$ws.Range("H122").Value = 70958.2
$ws.Range("I122").Value = 95323.09
$ws.Range("K122").Value = 285969.27
$ws.Range("M122").Value = -283519.27
$ws.Range("H126").Value = 3937.375
$ws.Range("I126").Value = 1715.1666
$ws.Range("J126").Value = 5270.7
$ws.Range("K126").Value = 5145.4998
$ws.Range("L126").Value = 15812.1
$ws.Range("M126").Value = -2675.4998
$ws.Range("N126").Value = -20752.1
$ws.Range("H132").Value = 4107.6875
$ws.Range("I132").Value = 1793.0454
$ws.Range("J132").Value = 9199.9
$ws.Range("K132").Value = 5379.1362
$ws.Range("L132").Value = 27599.7
$ws.Range("M132").Value = -2849.1362
$ws.Range("N132").Value = -32659.7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3485.9375
$ws.Range("I61").Value = 2454.5293
$ws.Range("K61").Value = 2454.5293
$ws.Range("M61").Value = -2252.5293
$ws.Range("H100").Value = 4600.857
$ws.Range("I100").Value = 3499.5
$ws.Range("K100").Value = 3499.5
$ws.Range("M100").Value = -2958.5
$ws.Range("H113").Value = 3485.9375
$ws.Range("I113").Value = 2454.5293
$ws.Range("K113").Value = 2454.5293
$ws.Range("M113").Value = -284.5293000000001
$ws.Range("H136").Value = 12201.883
$ws.Range("I136").Value = 4347.4287
$ws.Range("J136").Value = 17700
$ws.Range("K136").Value = 13042.2861
$ws.Range("L136").Value = 53100
$ws.Range("M136").Value = -10492.2861
$ws.Range("N136").Value = -58200
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2575
$ws.Range("I96").Value = 2400
$ws.Range("J96").Value = 2866.6667
$ws.Range("K96").Value = 2400
$ws.Range("L96").Value = 2866.6667
$ws.Range("M96").Value = -1027
$ws.Range("N96").Value = -5612.6667
$ws.Range("H106").Value = 45000
$ws.Range("J106").Value = 45000
$ws.Range("L106").Value = 45000
$ws.Range("N106").Value = -47524
$ws.Range("H132").Value = 25032772
$ws.Range("I132").Value = 45465820
$ws.Range("J132").Value = 59048
$ws.Range("K132").Value = 136397460
$ws.Range("L132").Value = 177144
$ws.Range("M132").Value = -136394930
$ws.Range("N132").Value = -182204
